# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2 through 91) from serial date 45192 (2023-09-23) to
# serial date 45202 (2023-10-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = 91

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 45192) {
        $cell.Value = 45202
    }
}
